# Update "CORRIDAS TK MULTIMARCAS.xlsx"
# - MES 02: add two new trip entries (29 Nov and 30 Nov 2023) before the
#   TOTAL row, update the TOTAL formulas/merge to match, and refresh the
#   page setup (fit to page, margins, zoom) + view scroll/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MES 02")

# --- Insert two rows above the TOTAL row (currently row 62) -------------
$ws.Rows.Item(62).Insert(-4121)
$ws.Rows.Item(62).Insert(-4121)

# Copy the visual formatting (borders/number formats) of row 59 - the
# existing template for a single-entry ("N deg"=1) boxed row - onto the
# two freshly inserted rows.
$ws.Range("A59:D59").Copy()
$ws.Range("A62:D62").PasteSpecial(-4122)
$ws.Range("A59:D59").Copy()
$ws.Range("A63:D63").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Fill in the new data -------------------------------------------------
$ws.Cells.Item(62, 1).Value = 1
$ws.Cells.Item(62, 2).Value = 45259
$ws.Cells.Item(62, 3).Value = "VILA ROSARIO"
$ws.Cells.Item(62, 4).Value = 10

$ws.Cells.Item(63, 1).Value = 1
$ws.Cells.Item(63, 2).Value = 45260
$ws.Cells.Item(63, 3).Value = "COSTA ESMERALDA"
$ws.Cells.Item(63, 4).Value = 15

# --- Fix up the TOTAL row (now row 64) ------------------------------------
$ws.Cells.Item(64, 1).Formula = "=SUM(A2:A63)"
$ws.Cells.Item(64, 4).Formula = "=SUM(D2:D63)"

# --- Page setup: fit to page, margins (cm), zoom --------------------------
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1
$ws.PageSetup.Zoom = 81

$ws.PageSetup.LeftMargin = $excel.CentimetersToPoints(3.6)
$ws.PageSetup.RightMargin = $excel.CentimetersToPoints(0.6)
$ws.PageSetup.TopMargin = $excel.CentimetersToPoints(0.5)
$ws.PageSetup.BottomMargin = $excel.CentimetersToPoints(0.4)
$ws.PageSetup.HeaderMargin = $excel.CentimetersToPoints(0.8)
$ws.PageSetup.FooterMargin = $excel.CentimetersToPoints(0.8)

# --- View: scroll position + active selection ------------------------------
$ws.Activate()
$ws.Range("I64").Select()
$excel.ActiveWindow.ScrollRow = 41
$excel.ActiveWindow.ScrollColumn = 1
